$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45203
}

$ws.Range("H2").Value = 21
$ws.Range("J2").Value = 19
$ws.Range("K2").Value = 3
$ws.Range("O2").Value = 22
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 48

$species = "Knärot`nRynkskinn`nTaggfingersvamp`nBarrviolspindling`nGarnlav`nGrantaggsvamp`nGul taggsvamp`nGultoppig fingersvamp`nJärpe`nMindre hackspett`nMotaggsvamp`nNordfladdermus`nOrange taggsvamp`nRosenticka`nRödvingetrast`nSkuggviol`nSpillkråka`nStjärntagging`nSvartvit taggsvamp`nTalltita`nTretåig hackspett`nUllticka`nBarkticka`nBronshjon`nFjällig taggsvamp s.str.`nHagfingersvamp`nKorallblylav`nPlattlummer`nRödgul trumpetsvamp`nSkinnlav`nSotriska`nSpindelblomster`nSvavelriska`nThomsons trägnagare`nTvåblad`nVedticka`nVågbandad barkbock`nÖgonpyrola`nStörre brunfladdermus`nVattenfladdermus`nÅkergroda`nVanlig groda`nFläcknycklar`nGrönvit nattviol`nNattviol`nBlåsippa`nLopplummer`nRevlummer"

$ws.Range("R2").Value = $species

# The sheet uses a fixed custom row height everywhere; undo Excel's
# automatic row-height growth triggered by the longer wrapped text.
$ws.Rows.Item(2).RowHeight = 15
